# Updated cryptos list - applies Price (D) and Volume(1h) (E) changes per row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.760.41"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").Value = "2.250.74"
$ws.Range("E3").Value = "  +0.98%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "112.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "294.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.47%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.631"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.02%  "

$ws.Range("E8").Value = "  -0.47%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.607"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0926"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.98%  "

$ws.Range("E13").Value = "  -0.97%  "

$ws.Range("E14").Value = "  +23.42%  "

$ws.Range("E15").Value = "  -0.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.40%  "

$ws.Range("D17").Value = "2.589.49"
$ws.Range("E17").Value = "  +0.79%  "

$ws.Range("D18").Value = "2.281.28"
$ws.Range("E18").Value = "  +1.94%  "

$ws.Range("D19").Value = "42.725.47"
$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("E20").Value = "  +0.46%  "

$ws.Range("E21").Value = "  +6.78%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +17.15%  "

$ws.Range("E24").Value = "  +5.72%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "255.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.97%  "

$ws.Range("E27").Value = "  -0.16%  "

$ws.Range("E28").Value = "  -2.76%  "

$ws.Range("E29").Value = "  -0.39%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.88%  "

$ws.Range("E33").Value = "  -3.34%  "

$ws.Range("E34").Value = "  +0.40%  "

$ws.Range("E35").Value = "  +2.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.57%  "

$ws.Range("E37").Value = "  +0.65%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0378"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.50%  "

$ws.Range("E40").Value = "  -0.97%  "

$ws.Range("E41").Value = "  -5.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.63%  "

$ws.Range("E43").Value = "  +0.27%  "

$ws.Range("E44").Value = "  -0.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.63%  "

$ws.Range("E46").Value = "  +0.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.70%  "

$ws.Range("E49").Value = "  +3.26%  "

$ws.Range("E50").Value = "  +2.68%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.33%  "
